# FindPattern section: the "Custom file format for specifying patterns."
# todo item has been completed, so remove it from the list entirely
# (commit message: "Remove completed item from todo list.").
$d = $word.ActiveDocument

$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "Custom file format for specifying patterns.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    $para = $searchRange.Paragraphs.Item(1)
    # Delete the whole paragraph, including its trailing paragraph mark,
    # so the list collapses cleanly onto the following "PeLib" heading.
    $wholePara = $d.Range($para.Range.Start, $para.Range.End + 1)
    $wholePara.Delete()
}

# Word stamps the "_GoBack" bookmark at the location of the most recent
# edit. Since that edit happened right after deleting the paragraph above
# (the cursor lands back near where typing/deleting last occurred, a few
# characters into "Pattern generator."), move the bookmark there.
$patternRange = $d.Content
$patternFound = $patternRange.Find.Execute(
    "Pattern generator.",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($patternFound) {
    $goBackPos = $patternRange.Start + 3
    $goBackRange = $d.Range($goBackPos, $goBackPos)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
